# "Loan RBI, Variable Instalments"
#
# The repayment schedule gains a new (blank) column so that a "Variable
# Instalments" / Late-tracking column can be inserted ahead of the existing
# "Late" column. Inserting a whole column shifts every cell from N onward
# one slot to the right (N->O, O->P, P->Q) and widens the used range from
# A1:P14 to A1:Q14 - exactly like a user choosing Insert > Sheet Columns
# from the ribbon while the "Late" column was selected.
#
# The author's last action before saving was on the "Repayment Schedule"
# sheet with cell L19 selected, so we finish by activating that sheet and
# selecting L19 (this is also what flips the workbook's stored activeTab
# and each sheet's tabSelected flag).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new blank column at N, pushing the old N/O/P columns to O/P/Q.
$ws.Columns("N:N").Insert()

# Make "Repayment Schedule" the active sheet/tab and restore the selection
# that was in place when the file was saved.
$ws.Activate()
[void]$ws.Range("L19").Select()
